$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (rows 13-16): reuse the existing "id text" style (same style as A5, A2, etc.) ---
$ws.Range("A5").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A16").PasteSpecial(-4122)

# --- New cell styles for column B (rows 13-16), matching two distinct fonts already used elsewhere ---
# Style for B13: font color FF010101 "Calibri  " (same font as cells like A2/A4/A6..)
$ws.Range("ZZ500").Font.Name = "Calibri  "
$ws.Range("ZZ500").Font.Color = 65793
$ws.Range("ZZ500").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("ZZ500").Clear()

# Style for B14:B16: theme color 1 "Calibri  " (same font as cells like A3/A9/A5)
$ws.Range("ZZ500").Font.Name = "Calibri  "
$ws.Range("ZZ500").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("ZZ500").Clear()

# --- Values (order matters, to reproduce the exact shared-string table order) ---
$ws.Range("A13").Value = "4255891924062617088"
$ws.Range("B13").Value = "XTE J1855-026"
$ws.Range("A14").Value = "4070968778552176768"
$ws.Range("B14").Value = "IGR J18027-2016"
$ws.Range("B15").Value = "EXO 1722-363"
$ws.Range("B16").Value = "OAO 1657-415"
$ws.Range("A16").Value = "5966391894137386240"
$ws.Range("A15").Value = "0000000000000000000"

# --- Column B width needs to grow to fit the new, longer star names ---
$ws.Columns("B").AutoFit()

# --- Selection moved, as in the saved workbook view state ---
$ws.Range("F19").Select()
